# Generate Report for Handback
#
# Fills in the "4624e03b-42b9-416c-bf73-a169687bcbbe" row (row 7) of the
# per-locale handback-status sheets (zh-cn, de-de) now that a handback was
# processed for that file: the "Latest Target File" becomes a hyperlink to
# the handed-back markdown, "Latest Handback File" / "Latest Handback
# DateTime" get filled in, and "Error Detail" records that the handback
# was based on a stale source revision.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # OLE/BGR encoding of RGB(100,149,237) = FF6495ED, matches the existing "HyperLink" cell style

function Set-HandbackRow7 {
    param(
        [string]$SheetName,
        [string]$TargetFile,     # Latest Target File (column I) display text
        [string]$TargetUrl,      # Latest Target File (column I) hyperlink URL
        [string]$HandbackFile,   # Latest Handback File (column J)
        [string]$HandbackTime,   # Latest Handback DateTime (column K)
        [string]$ErrorDetail     # Error Detail (column P)
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $iCell = $ws.Range("I7")
    $ws.Hyperlinks.Add($iCell, $TargetUrl, "", "", $TargetFile)
    $iCell.Font.Underline = $true
    $iCell.Font.Color = $hyperlinkColor

    $ws.Range("J7").Value = $HandbackFile
    $ws.Range("K7").Value = $HandbackTime
    $ws.Range("P7").Value = $ErrorDetail
}

$latestUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a6f2cc4ae622b440917870474dd26e8761b3e62/e2e/4624e03b-42b9-416c-bf73-a169687bcbbe.md"

$currentUrlZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/63516e4b40c3a38ab65f986b764968fd786b8249/e2e/4624e03b-42b9-416c-bf73-a169687bcbbe.md"
$errorDetailZhCn = "The version of handback file is not the latest, current: $currentUrlZhCn, latest: $latestUrl."

Set-HandbackRow7 "zh-cn" `
    "4624e03b-42b9-416c-bf73-a169687bcbbe.md" `
    $currentUrlZhCn `
    "4624e03b-42b9-416c-bf73-a169687bcbbe.a0da62b9bf8155d231e5f00fdde60240f9dc9e80.zh-cn.xlf" `
    "2016-08-30 21:02:58" `
    $errorDetailZhCn

$currentUrlDeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/63516e4b40c3a38ab65f986b764968fd786b8249/e2e/4624e03b-42b9-416c-bf73-a169687bcbbe.md"
$errorDetailDeDe = "The version of handback file is not the latest, current: $currentUrlDeDe, latest: $latestUrl."

Set-HandbackRow7 "de-de" `
    "4624e03b-42b9-416c-bf73-a169687bcbbe.md" `
    $currentUrlDeDe `
    "4624e03b-42b9-416c-bf73-a169687bcbbe.a0da62b9bf8155d231e5f00fdde60240f9dc9e80.de-de.xlf" `
    "2016-08-30 21:03:13" `
    $errorDetailDeDe
